# Insert a new Mango price-record row before the current row 191 (Ecuador,
# week of 2021-11-29), shifting every subsequent record down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("191:191").Insert()

$ws.Range("A191").Value = 9
$ws.Range("B191").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C191").Value = "Metropolitana"
$ws.Range("D191").Value = 44529
$ws.Range("E191").Value = 13
$ws.Range("F191").Value = "Fruta"
$ws.Range("G191").Value = 100108
$ws.Range("H191").Value = "Tropicales y subtropicales"
$ws.Range("I191").Value = 100108002
$ws.Range("J191").Value = "Mango"
$ws.Range("K191").Value = "Sin especificar"
$ws.Range("L191").Value = "Primera"
$ws.Range("M191").Value = 380
$ws.Range("N191").Value = 6000
$ws.Range("O191").Value = 6500
$ws.Range("P191").Value = 6263
$ws.Range("Q191").Value = "$/bandeja 4 kilos"
$ws.Range("R191").Value = "Ecuador"
$ws.Range("S191").Value = 1566
$ws.Range("T191").Value = 4
